$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f = $sec.Footers(1)
$rng = $f.Range
$rng.Fields.Add($rng, 33)
Write-Host "done"
